$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.705
$ws.Range("C9").Value = -10.85
$ws.Range("E12").Value = 17.739
$ws.Range("E14").Value = 16.70500000000001
$ws.Range("C18").Value = -11.62
$ws.Range("C20").Value = -12.358
$ws.Range("E26").Value = 16.525
$ws.Range("C27").Value = -13.224
$ws.Range("E27").Value = 16.366
$ws.Range("E29").Value = 16.941
$ws.Range("C35").Value = -12.197
$ws.Range("E37").Value = 16.852
$ws.Range("E38").Value = 16.745
$ws.Range("E51").Value = 16.743
$ws.Range("E52").Value = 16.83
$ws.Range("E55").Value = 16.494
$ws.Range("C69").Value = -10.672
$ws.Range("E69").Value = 17.396
$ws.Range("E70").Value = 17.629
$ws.Range("C76").Value = -13.201
$ws.Range("C78").Value = -12.5
$ws.Range("E81").Value = 16.434
$ws.Range("C82").Value = -11.991
$ws.Range("C83").Value = -13.192
$ws.Range("E83").Value = 16.731
$ws.Range("C93").Value = -11.642
$ws.Range("E102").Value = 16.784
